# 3dunet session annotations.xlsx - prepare train3dunet 1.8.2 chpt-231225-3
#
# Context: sheet1 is a session log. Row 109 was a blank "TBD" template row
# (all cells pre-filled with the shared "TBD" placeholder string plus a
# handful of default values/formulas). This edit:
#   1. fills in results for the two earlier sessions (rows 107 & 108) that
#      turned out to fail for a gcc-related reason,
#   2. turns the former blank template row 109 into a real log entry
#      (session 231225-2),
#   3. turns row 110 (another near-blank placeholder row) into a new log
#      entry (session 231225-3) carrying the same formula/style template
#      row 109 used to have, and
#   4. appends a fresh blank template row 111 (same template again) so the
#      sheet keeps two spare rows at the bottom, as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Existing rows 107 & 108: record the (failed, gcc-related) outcome.
#    New shared strings are minted in this exact order so they land at
#    sharedStrings indices 469-474, matching the target workbook.
# ---------------------------------------------------------------------

# -- new row 109 becomes session "231225-2" --
$ws.Range("A109").Value = "231225-2"                                             # -> new string 469

# -- go back and annotate rows 107/108 with the gcc-related failure reason --
$ws.Range("BD107").Value = "gcc related"                                         # -> new string 470
$ws.Range("BD108").Value = "gcc related"                                         # reuse 470

# -- finish describing the new 231225-2 entry --
$ws.Range("E109").Value = "Copy test_config-231225-0.yml and rerun predict3dunet on it from an interactive A100 session."   # -> new string 471

# -- record the (shared) failure result text on rows 107/108 --
$ws.Range("G107").Value = "Fail. Gcc related error."                             # -> new string 472
$ws.Range("G108").Value = "Fail. Gcc related error."                             # reuse 472

# -- new row 110 becomes session "231225-3" --
$ws.Range("A110").Value = "231225-3"                                             # -> new string 473
$ws.Range("E110").Value = "Copy train_config-231225-1.yml and rerun predict3dunet on it from an interactive A100 session."  # -> new string 474

# ---------------------------------------------------------------------
# 2) Row 107 / 108 numeric + boolean-ish cells.
# ---------------------------------------------------------------------
$ws.Range("H107").Value = 0
$ws.Range("BC107").Value = 1

$ws.Range("F108").Value = "Success (no error)"
$ws.Range("H108").Value = 0
$ws.Range("BC108").Value = 1

# ---------------------------------------------------------------------
# 3) Flesh out row 109 (session 231225-2) - the rest of its columns.
# ---------------------------------------------------------------------
$ws.Range("B109").Value = "predict3dunet 1.6.0"
$ws.Range("C109").Value = "ResidualUNet3D"
$ws.Range("D109").Value = "3DUnet_lightsheet_boundary"
$ws.Range("F109").Value = "Success (no error)"
$ws.Range("G109").Value = "Fail. Patch shape invalid error."
$ws.Range("H109").Value = 0
$ws.Range("BC109").Value = 1
$ws.Range("BD109").Value = "ValueError: requested an output size of torch.Size([13, 121, 33]), but valid sizes range from [11, 119, 31] to [12, 120, 32] (for an input of torch.Size([6, 60, 16]))"
$ws.Range("BE109").Value = "  File ""/home/dwalth/data/conda/envs/3dunet/lib/python3.11/site-packages/torch/nn/modules/conv.py"", line 662, in _output_padding`n    raise ValueError(("
$ws.Range("BE109").WrapText = $true

# Row 109 grew a second line of text -> Excel bumped its row height.
$ws.Rows.Item(109).RowHeight = 28.8

# ---------------------------------------------------------------------
# 4) Build out row 110 as a full copy of the old "TBD" template (same
#    shape row 109 used to have), then overwrite it with the new
#    231225-3 entry's known data. Columns not yet known stay "TBD".
# ---------------------------------------------------------------------
$templateCols = @("B","C","D","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
  "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW", `
  "BA","BB","BC","BD","BE")

foreach ($row in 110,111) {
    foreach ($col in @("A","E","F","G") + $templateCols) {
        $cell = $ws.Range($col + $row)
        if ($cell.Value2 -eq $null -or $cell.Value2 -eq "") {
            $cell.Value = "TBD"
        }
    }
    $ws.Range("AK" + $row).Value = "TBD: formula is TBD"
    $ws.Range("AL" + $row).Value = "NVIDIA A100-SXM4-80GB"
    $ws.Range("AS" + $row).Value = "NA"
    $ws.Range("AW" + $row).Value = "NA"
    $ws.Range("AX" + $row).Formula = "=AP" + $row + "-AT" + $row
    $ws.Range("AY" + $row).Formula = "=AQ" + $row + "-AU" + $row
    $ws.Range("AZ" + $row).Formula = "=AR" + $row + "-AV" + $row
}

# Now that the template has been stamped onto row 110, overwrite it with
# the actual 231225-3 entry's known values (G110 stays "TBD" - outcome
# not yet known).
$ws.Range("B110").Value = "train3dunet 1.8.2"
$ws.Range("C110").Value = "ResidualUNet3D"
$ws.Range("D110").Value = "3DUnet_lightsheet_boundary"
$ws.Range("F110").Value = "Success (no error)"
